# Apply the commit: "Removed Birth age program from all books,
# removed IYCF 2&3 from default books"
#
# - Sheets "Programs to include", "Coverage scenario", "Budget scenario":
#   remove the rows for "Birth age program", "IYCF 2" and "IYCF 3".
# - Sheet "Coverage scenario": "Micronutrient powders" and
#   "Zinc for treatment + ORS" now have a default coverage of 0.95
#   (column D) instead of being blank.
# - Sheet "Optimization options": the objectives cell no longer
#   includes "child_deaths" and the "additional funds" amount drops
#   from 10,000,000 to 5,000,000.

$wb = $excel.ActiveWorkbook

$namesToRemove = @("Birth age program", "IYCF 2", "IYCF 3")

function Remove-NamedRows($ws, $names) {
    $ur = $ws.UsedRange
    $lastRow = $ur.Rows.Count
    for ($r = $lastRow; $r -ge 1; $r--) {
        $val = $ws.Cells.Item($r, 1).Value2
        if ($names -contains $val) {
            $ws.Rows.Item($r).Delete()
        }
    }
}

# --- Sheet 1: "Programs to include" ---
$wsPrograms = $wb.Worksheets.Item("Programs to include")
Remove-NamedRows $wsPrograms $namesToRemove

# --- Sheet 2: "Coverage scenario" ---
$wsCoverage = $wb.Worksheets.Item("Coverage scenario")
Remove-NamedRows $wsCoverage $namesToRemove

$ur = $wsCoverage.UsedRange
$lastRow = $ur.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $name = $wsCoverage.Cells.Item($r, 1).Value2
    if ($name -eq "Micronutrient powders" -or $name -eq "Zinc for treatment + ORS") {
        $wsCoverage.Cells.Item($r, 4).Value = 0.95
    }
}

# --- Sheet 3: "Budget scenario" ---
$wsBudget = $wb.Worksheets.Item("Budget scenario")
Remove-NamedRows $wsBudget $namesToRemove

# --- Sheet 4: "Optimization options" ---
$wsOpt = $wb.Worksheets.Item("Optimization options")
$wsOpt.Range("C2").Value = "thrive"
$wsOpt.Range("E2").Value = 5000000
